# Upgrade the "SNMP" test-case sheet:
#  - authPro (B2) changes from "MD5|SHA" to "MD5"
#  - privPro (D2) changes from "DES|AES128" to "DES"
#  - the sheet's active selection moves to D2 (and the old horizontal
#    scroll position pinned to column B is cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SNMP")
$ws.Activate()

$ws.Range("B2").Value = "MD5"
$ws.Range("D2").Value = "DES"

$ws.Range("D2").Select()
